$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Server-side status refresh: HK1032 / Yeshi / Tandoori Chicken Biriyani -> Delivered ---
$ws.Cells.Item(48, 8).Value = "Delivered"

# --- Append four new order rows (HK1034-HK1037), all placed/dated 12/16/2025 ---
$newRows = @(
    @("HK1034", "12/16/2025", "Sunitha Woodlands", "Dosa Batter", 7, 2, 14, "Accepted"),
    @("HK1035", "12/16/2025", "Pranav", "Daily Subscription", 8, 1, 8, "Accepted"),
    @("HK1036", "12/16/2025", "Shalini Raju", "Veg Curries", 10, 2, 20, "Accepted"),
    @("HK1037", "12/16/2025", "Hemanth", "Boiled egg curry with Jeera rice", 15, 1, 15, "Accepted")
)

$firstRow = 50
$lastRow = $firstRow + $newRows.Count - 1

# Keep the Date column (B) as literal text ("12/16/2025"), matching the rest of the
# sheet's text-formatted dates, instead of letting Excel auto-convert it to a date serial.
$dateRange = $ws.Range("B$firstRow" + ":B$lastRow")
$dateRange.NumberFormat = "@"

$r = $firstRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Reset the Date column's style back to the default (Normal) now that the text value
# is locked in, so no stray number-format styling is left applied to the cells.
$dateRange.Style = "Normal"
